$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Duplicate the "iter_04" sheet (inherits column widths / number
#    formats / base styling) and place the copy right after it, then
#    rename it to "iter_05".
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("iter_04")
$src.Copy($null, $src)
$new = $wb.Worksheets.Item($src.Index + 1)
$new.Name = "iter_05"

# ---------------------------------------------------------------------
# 2) Update the calibrated parameter values on iter_05. Rows/cols match
#    the calibration-log layout used by iter_02 / iter_04:
#      col B..G = local_bus, express_bus, ferry, light_rail, heavy_rail,
#                 commuter_rail ; col H = walk
# ---------------------------------------------------------------------
$values = @{
    "B3"=3;  "C3"=3;  "D3"=1;  "E3"=3;  "F3"=3;  "G3"=1;
    "B4"=3;  "C4"=3;  "D4"=1;  "E4"=3;  "F4"=3;  "G4"=1;
    "G5"=0.2;
    "H8"=2;
}
foreach ($addr in $values.Keys) {
    $new.Range($addr).Value = $values[$addr]
}

# ---------------------------------------------------------------------
# 3) Re-apply the "changed since previous iteration" bold-red highlight.
#    Cells whose value differs from the prior iteration (iter_04) are
#    bold + red; everything else reverts to plain formatting. Row 6
#    (transfer_wait_perception_factor) no longer changed, so it loses
#    its old highlight.
#    NB: plain (unhighlighted) cells are restyled *before* the
#    bold/red ones - doing it in the opposite order confuses the
#    engine's style de-duplication (a plain numeric cell can end up
#    reusing a stale quote-prefixed style) once a neighbouring
#    quote-prefixed text cell has already been touched.
# ---------------------------------------------------------------------
$highlighted = @("B3","C3","E3","F3","G3","B4","C4","D4","E4","F4","G4","G5","H8")
$plain       = @("B2","C2","D2","E2","F2","G2","D3","B5","C5","D5","E5","F5","B6","C6","D6","E6","F6","G6")

foreach ($addr in $plain) {
    $cell = $new.Range($addr)
    $cell.Font.Bold = $false
    $cell.Font.Color = 0
}
foreach ($addr in $highlighted) {
    $cell = $new.Range($addr)
    $cell.Font.Bold = $true
    $cell.Font.Color = 255
}

# ---------------------------------------------------------------------
# 4) Sheet view cosmetics for all three sheets, per the recorded diff.
# ---------------------------------------------------------------------
$iter02 = $wb.Worksheets.Item("iter_02")
$iter02.Activate()
$iter02.Range("A7").Select() | Out-Null
$excel.ActiveWindow.Zoom = 150

$iter04 = $wb.Worksheets.Item("iter_04")
$iter04.Activate()
$iter04.Range("H8").Select() | Out-Null
$excel.ActiveWindow.Zoom = 150

$new.Activate()
$new.Range("C8").Select() | Out-Null
$excel.ActiveWindow.Zoom = 150

# Return focus to iter_02, which stays the selected/active tab.
$iter02.Activate()
